# Commit: "regen save_data to use K instead of Strike#, regen std/mean,
# calc and write s_vals". For this sheet, the regenerated pipeline produced
# new values for column G ("K") for each game row (rows 2-76). The sheet
# stores raw data (no formulas), so we just write the recalculated values
# directly into column G, matching the upstream regen.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row number (A1-style row) -> new K value
$kValues = @{
    2 = 2;
    3 = 2;
    4 = 2;
    5 = 2;
    6 = 3;
    7 = 2;
    8 = 1;
    9 = 1;
    10 = 1;
    11 = 2;
    12 = 3;
    13 = 2;
    14 = 0;
    15 = 1;
    16 = 1;
    17 = 2;
    18 = 0;
    19 = 2;
    20 = 1;
    21 = 1;
    22 = 3;
    23 = 3;
    24 = 0;
    25 = 3;
    26 = 4;
    27 = 2;
    28 = 2;
    29 = 0;
    30 = 1;
    31 = 3;
    32 = 2;
    33 = 2;
    34 = 3;
    35 = 3;
    36 = 4;
    37 = 1;
    38 = 1;
    39 = 1;
    40 = 1;
    41 = 1;
    42 = 3;
    43 = 2;
    44 = 1;
    45 = 1;
    46 = 2;
    47 = 0;
    48 = 2;
    49 = 0;
    50 = 2;
    51 = 0;
    52 = 5;
    53 = 3;
    54 = 0;
    55 = 2;
    56 = 0;
    57 = 0;
    58 = 2;
    59 = 3;
    60 = 2;
    61 = 0;
    62 = 2;
    63 = 1;
    64 = 2;
    65 = 5;
    67 = 2;
    68 = 3;
    69 = 2;
    70 = 1;
    71 = 2;
    72 = 1;
    73 = 3;
    74 = 1;
    75 = 2;
    76 = 1;
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item([int]$row, 7).Value = $kValues[$row]
}
